$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $val) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

# Row 2
Set-TextValue "D2" "62.458.52"
Set-TextValue "E2" "  +5.03%  "

# Row 3
Set-TextValue "D3" "3.324.53"
Set-TextValue "E3" "  +4.69%  "

# Row 4
Set-TextValue "E4" "  +0.05%  "

# Row 5
Set-TextValue "D5" "552.96"
Set-TextValue "E5" "  +3.36%  "

# Row 6
Set-TextValue "D6" "151.01"
Set-TextValue "E6" "  +4.69%  "

# Row 7
Set-TextValue "D7" "1.00"
Set-TextValue "E7" "  +0.04%  "

# Row 8
Set-TextValue "D8" "0.528"
Set-TextValue "E8" "  +1.34%  "

# Row 9
Set-TextValue "D9" "7.51"
Set-TextValue "E9" "  +3.16%  "

# Row 10
Set-TextValue "E10" "  +4.12%  "

# Row 11
Set-TextValue "D11" "0.437"
Set-TextValue "E11" "  +1.44%  "

# Row 12
Set-TextValue "D12" "3.898.88"
Set-TextValue "E12" "  +4.69%  "

# Row 13
Set-TextValue "E13" "  -1.00%  "

# Row 14
Set-TextValue "E14" "  +4.97%  "

# Row 15
Set-TextValue "D15" "26.94"
Set-TextValue "E15" "  +3.47%  "

# Row 16
Set-TextValue "D16" "62.418.02"
Set-TextValue "E16" "  +4.92%  "

# Row 17
Set-TextValue "D17" "3.329.51"
Set-TextValue "E17" "  +5.87%  "

# Row 18
Set-TextValue "D18" "6.50"
Set-TextValue "E18" "  +5.04%  "

# Row 19
Set-TextValue "D19" "13.80"
Set-TextValue "E19" "  +6.67%  "

# Row 21
Set-TextValue "D21" "383.82"
Set-TextValue "E21" "  +1.83%  "

# Row 22
Set-TextValue "E22" "  +0.20%  "

# Row 23
Set-TextValue "D23" "0.537"
Set-TextValue "E23" "  +1.77%  "

# Row 24
Set-TextValue "D24" "70.84"
Set-TextValue "E24" "  +1.15%  "

# Row 25
Set-TextValue "E25" "  +4.15%  "

# Row 26
Set-TextValue "D26" "8.83"
Set-TextValue "E26" "  +0.36%  "

# Row 27
Set-TextValue "D27" "0.0₃0967"
Set-TextValue "E27" "  +7.59%  "

# Row 28
Set-TextValue "E28" "  +0.19%  "

# Row 29
Set-TextValue "E29" "  +3.52%  "

# Row 30
Set-TextValue "E30" "  +3.69%  "

# Row 31
Set-TextValue "B31" "EthereumClassic"
Set-TextValue "C31" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D31" "22.96"
Set-TextValue "E31" "  +2.97%  "

# Row 32
Set-TextValue "B32" "Fetch.AI"
Set-TextValue "C32" "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D32" "1.31"
Set-TextValue "E32" "  +10.94%  "

# Row 33
Set-TextValue "E33" "  +3.95%  "

# Row 34
Set-TextValue "D34" "6.72"
Set-TextValue "E34" "  +3.94%  "

# Row 35
Set-TextValue "E35" "  +10.98%  "

# Row 36
Set-TextValue "D36" "159.53"
Set-TextValue "E36" "  +2.02%  "

# Row 37
Set-TextValue "D37" "1.87"
Set-TextValue "E37" "  +11.60%  "

# Row 38
Set-TextValue "D38" "26.94"
Set-TextValue "E38" "  +5.82%  "

# Row 39
Set-TextValue "D39" "2.845.53"
Set-TextValue "E39" "  +4.09%  "

# Row 40
Set-TextValue "E40" "  +2.67%  "

# Row 41
Set-TextValue "D41" "0.0314"
Set-TextValue "E41" "  +8.23%  "

# Row 42
Set-TextValue "D42" "4.33"
Set-TextValue "E42" "  +1.37%  "

# Row 43
Set-TextValue "E43" "  +3.48%  "

# Row 44
Set-TextValue "D44" "40.53"
Set-TextValue "E44" "  +3.03%  "

# Row 45
Set-TextValue "E45" "  +4.02%  "

# Row 46
Set-TextValue "E46" "  +7.17%  "

# Row 47
Set-TextValue "D47" "3.370.17"
Set-TextValue "E47" "  +4.72%  "

# Row 48
Set-TextValue "D48" "0.104"
Set-TextValue "E48" "  +4.12%  "

# Row 49
Set-TextValue "D49" "6.30"
Set-TextValue "E49" "  +1.95%  "

# Row 50
Set-TextValue "D50" "0.808"
Set-TextValue "E50" "  +5.15%  "

# Row 51
Set-TextValue "D51" "281.29"
Set-TextValue "E51" "  +7.50%  "
